$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The applicant has extensive experience in working on projects involving ReactJS, Django, Python, JavaScript, Computer Vision, Image Processing, Data Analysis, and PyTorch, which are all relevant skills required for the Frontend Engineer Intern position. The applicant has demonstrated the ability to align web design and user experience requirements, optimize web pages, and maintain brand consistency through projects like ''LLMGuard'' and ''Literature Society IITJ Website''. Additionally, the applicant has shown proficiency in collaborating with back-end developers, graphic designers, and UX designers, as mentioned in the job description. However, the applicant lacks direct experience with CSS and NextJS, which are required skills for the position. Despite this, the applicant''s diverse technical skills and experience make them a strong candidate for the role.'
$ws.Range("E2").Value = 'The applicant has expressed a strong interest in working in Japan, citing admiration for the country''s work culture, a willingness to learn a new language, and a desire to be part of cutting-edge AI/ML technologies. The applicant also demonstrates an adaptable nature and emphasizes the importance of teamwork, problem-solving, and communication. While expressing attachment to their Indian roots, the applicant''s openness to adapting to new cultures and their clear expectations from the company showcase a positive and determined attitude towards working in Japan.'
$ws.Range("C3").Value = 'The applicant has a diverse set of technical skills including but not limited to Docker, NodeJS, ExpressJS, Firebase, and ReactJS. Although the applicant''s projects showcase proficiency in developing web applications using relevant technologies, it is evident that the applicant lacks experience with NextJS, which is a specific requirement mentioned in the job description. Furthermore, while the projects demonstrate skills in CSS, JavaScript, and frontend development, there is no explicit mention of experience with ReactJS, a crucial requirement for the Frontend Engineer Intern position. Thus, while the applicant possesses a strong technical foundation, there are notable gaps in skills directly aligned with the company''s requirements.'
$ws.Range("E3").Value = 'The applicant has exhibited strong adaptability, communication, and problem-solving skills through the provided answers. The willingness to join a foreign work culture, especially in Japan, reflects open-mindedness and willingness to adapt to new environments. Furthermore, the applicant has expressed an interest in learning a new language while working in Japan, indicating a genuine interest in the culture and professional growth opportunities that Japan offers. The responses also demonstrate a positive attitude towards teamwork and collaboration. However, the applicant''s stated preference to not work alone may pose a minor challenge, hence the score of 4 is appropriate.'
$ws.Range("C4").Value = 'The applicant has a strong background in Python, PyTorch, Tensorflow, and Django Rest Framework, which are highly relevant to frontend engineering. While the applicant has not directly mentioned NextJS, they have experience in ReactJS, which can be considered transferable. Additionally, the applicant has showcased their expertise in computer vision, image processing, and machine learning, which demonstrates a strong technical foundation. Their projects exhibit proficiency in ReactJS, Django Rest Framework and demonstrating their ability to collaborate with back-end developers, aligning with the requirements of the job description. However, the applicant lacks direct mention of JavaScript and CSS, which are essential for frontend development. Despite this, their strong skills in other areas make them a potential candidate for the role, showcasing the ability to learn and adapt to new technologies, thus justifying a high score.'
$ws.Range("E4").Value = 'The applicant has shown a willingness to adapt to new cultures and work in Japan, appreciating the work culture and expressing interest in learning a new language. Moreover, they have emphasized their commitment to contributing to cutting-edge AI/ML technologies and learning new skills, demonstrating a proactive and enthusiastic attitude. The applicant''s experience in collaborating within a team, their ability to grasp new concepts, and their openness to improvements highlight their strong potential for successful integration into the company''s work environment in Japan, hence the score of 4.'
